$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 187
$ws.Range("I2").Value = 458
$ws.Range("J2").Value = 1999
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 583
$ws.Range("M2").Value = 32
$ws.Range("N2").Value = 358
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 24
$ws.Range("S2").Value = 186
$ws.Range("T2").Value = 352
$ws.Range("U2").Value = 30
$ws.Range("V2").Value = 2965
$ws.Range("X2").Value = 3147
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 35
$ws.Range("AA2").Value = 19
